$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.859.82'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.840.26'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('D4').Value = '''1.006'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '''309.02'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').Value = '''1.004'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.4698'
$ws.Range('E7').Value = '  +3.49%  '
$ws.Range('D8').Value = '''0.3657'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('D9').Value = '''0.07140'
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').Value = '''0.9225'
$ws.Range('E10').Value = '  +3.63%  '
$ws.Range('D11').Value = '''19.53'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '''0.07673'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = '1.892.76'
$ws.Range('E13').Value = '  +4.03%  '
$ws.Range('D14').Value = '''5.284'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '''6.386'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').Value = '''88.17'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '''1.007'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '''0.000008628'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = '''1.004'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '26.886.75'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '''14.44'
$ws.Range('E21').Value = '  +2.26%  '
$ws.Range('D22').Value = '''5.004'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D23').Value = '''10.58'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').Value = '''1.922'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '''151.69'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('D27').Value = '''2.005'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').Value = '''114.02'
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '''4.876'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').Value = '''0.08811'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('E31').Value = '  +2.54%  '
$ws.Range('D32').Value = '''1.174'
$ws.Range('E32').Value = '  +5.78%  '
$ws.Range('D33').Value = '''0.7450'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').Value = '''2.776'
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').Value = '''4.474'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').Value = '''1.086'
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').Value = '''0.01936'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.957'
$ws.Range('E38').Value = '  +1.60%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.05199'
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').Value = '''0.5187'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').Value = '''6.957'
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('D42').Value = '''0.1507'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '''8.138'
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('D44').Value = '''10.40'
$ws.Range('E44').Value = '  +4.17%  '
$ws.Range('D45').Value = '''0.4684'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = '''1.005'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').Value = '''101.21'
$ws.Range('E47').Value = '  +2.08%  '
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').Value = '''65.30'
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('D50').Value = '''0.06033'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '''0.8907'
$ws.Range('E51').Value = '  +5.31%  '
